$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Ticket Link" (column E) for the rows that got a Jira ticket id ---
# (set column E first, then column D, to match the shared-strings ordering
#  produced by the original edit)
$ws.Range("E2").Value = "EN-54135"
$ws.Range("E3").Value = "EN-54134"
$ws.Range("E7").Value = "EN-54131"
$ws.Range("E8").Value = "EN-54121"

# --- Update "Description" (column D) notes that were rewritten ---
$ws.Range("D2").Value = "Hi, I know this is hard "
$ws.Range("D3").Value = "Hi, OK, we should work on the EN-54134"
$ws.Range("D4").Value = "Hi, COOPPPPPPP"
$ws.Range("D9").Value = "Hi, Customers are required to wear a face covering in all of our stores. "

# --- Widen column D (Description) and give column E (Ticket Link) an explicit width ---
$ws.Columns.Item(4).ColumnWidth = 100.5
$ws.Columns.Item(5).ColumnWidth = 25.333333333333332

# --- Leave the selection where the user last clicked ---
[void]$ws.Range("D14").Select()
